# Update the "想去人数" (want-to-go count) figures that were refreshed by the
# gh-pages data regeneration (commit "Update gh-pages to output generated at 456a3b4").
# Column F on both the "展览" and "全部类型" sheets needs the following bumps.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"   = @{ 11 = 176; 13 = 6083; 19 = 474; 20 = 9153; 22 = 2459; 43 = 1530; 44 = 2515 }
    "全部类型" = @{ 15 = 176; 16 = 6083; 21 = 474; 22 = 9153; 24 = 2459; 41 = 1530; 42 = 2515 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item([int]$row, 6).Value = $rows[$row]
    }
}
